$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 279 (this shifts rows 279..347 down to 280..348,
# and copies the formatting, e.g. the date number format, from the row below).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new record.
$ws.Cells.Item(279, 1).Value = 5
$ws.Cells.Item(279, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(279, 3).Value = "Maule"
$ws.Cells.Item(279, 4).Value = Get-Date -Year 2022 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(279, 5).Value = 7
$ws.Cells.Item(279, 6).Value = 100114014
$ws.Cells.Item(279, 7).Value = "Betarraga"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 4000
$ws.Cells.Item(279, 11).Value = 700
$ws.Cells.Item(279, 12).Value = 700
$ws.Cells.Item(279, 13).Value = 700
$ws.Cells.Item(279, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(279, 15).Value = "Región del Maule"
$ws.Cells.Item(279, 16).Value = 140
$ws.Cells.Item(279, 17).Value = 5
$ws.Cells.Item(279, 18).Value = "Hortaliza"
